$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 12 (pushes settings..opt_hard rows down by one)
$ws.Rows.Item(12).Insert()

# Populate the new row with the "too many requests" translation strings
$ws.Cells.Item(12, 1).Value = "too_many_requests"
$ws.Cells.Item(12, 2).Value = "Please wait a minute before trying again."
$ws.Cells.Item(12, 3).Value = "Bitte warte eine Minute, bevor du es erneut versuchst."

# Re-apply the "code" column styling (bold italic 12pt, right aligned, centered vertically)
# used by every other row in column A
$ws.Cells.Item(12, 1).Font.Bold = $true
$ws.Cells.Item(12, 1).Font.Italic = $true
$ws.Cells.Item(12, 1).Font.Size = 12
$ws.Cells.Item(12, 1).HorizontalAlignment = -4152
$ws.Cells.Item(12, 1).VerticalAlignment = -4108

# Match the row height used throughout the rest of the table
$ws.Rows.Item(12).RowHeight = 15.75

# Widen columns B and C so the new, longer English/German strings fit
$ws.Columns.Item(2).ColumnWidth = 38
$ws.Columns.Item(3).ColumnWidth = 50.2

# Leave the cursor where it ended up after the edit
$ws.Range("C22").Select()
